$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells keep their text representation
# exactly as scraped (e.g. "0.9994", "26.460.85"), instead of Excel
# auto-converting them to numeric values which would change formatting
# and precision (e.g. "0.05330" -> 0.0533).

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.460.85"
$ws.Cells.Item(2, 5).Value = "  +2.72%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.730.23"
$ws.Cells.Item(3, 5).Value = "  +3.12%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9994"
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "244.02"
$ws.Cells.Item(5, 5).Value = "  +2.85%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9998"
$ws.Cells.Item(6, 5).Value = "  +0.05%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4798"
$ws.Cells.Item(7, 5).Value = "  +3.88%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2669"
$ws.Cells.Item(8, 5).Value = "  +3.05%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06232"
$ws.Cells.Item(9, 5).Value = "  +1.44%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.731.06"
$ws.Cells.Item(10, 5).Value = "  +3.19%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07118"
$ws.Cells.Item(11, 5).Value = "  +1.75%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.74"
$ws.Cells.Item(12, 5).Value = "  +5.62%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.6181"
$ws.Cells.Item(13, 5).Value = "  +6.77%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.546"
$ws.Cells.Item(14, 5).Value = "  +4.46%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "77.01"
$ws.Cells.Item(15, 5).Value = "  +2.19%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.9999"
$ws.Cells.Item(16, 5).Value = "  +0.12%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "26.467.48"
$ws.Cells.Item(17, 5).Value = "  +2.69%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "1.0000"
$ws.Cells.Item(18, 5).Value = "  +0.07%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000006936"
$ws.Cells.Item(19, 5).Value = "  +3.57%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.73"
$ws.Cells.Item(20, 5).Value = "  +2.59%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.955.04"
$ws.Cells.Item(21, 5).Value = "  +3.82%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.562"
$ws.Cells.Item(22, 5).Value = "  +2.08%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.899"
$ws.Cells.Item(23, 5).Value = "  +2.78%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.318"
$ws.Cells.Item(24, 5).Value = "  +1.53%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "136.47"
$ws.Cells.Item(25, 5).Value = "  +1.68%  "

$ws.Cells.Item(26, 5).Value = "  +2.36%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.790"
$ws.Cells.Item(27, 5).Value = "  +3.85%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.406"
$ws.Cells.Item(28, 5).Value = "  +1.62%  "

$ws.Cells.Item(29, 5).Value = "  +1.79%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.985"
$ws.Cells.Item(30, 5).Value = "  +0.89%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.07987"
$ws.Cells.Item(31, 5).Value = "  +4.03%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.736"
$ws.Cells.Item(32, 5).Value = "  +3.33%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04558"
$ws.Cells.Item(33, 5).Value = "  +4.86%  "

$ws.Cells.Item(34, 2).Value = "Frax"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.9996"
$ws.Cells.Item(34, 5).Value = "  +0.17%  "

$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.615"
$ws.Cells.Item(35, 5).Value = "  +0.47%  "

$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6404"
$ws.Cells.Item(36, 5).Value = "  +4.65%  "

$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9900"
$ws.Cells.Item(37, 5).Value = "  +3.94%  "

$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.9448"
$ws.Cells.Item(38, 5).Value = "  +1.27%  "

$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.984"
$ws.Cells.Item(39, 5).Value = "  +6.14%  "

$ws.Cells.Item(40, 2).Value = "Quant"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "107.50"
$ws.Cells.Item(40, 5).Value = "  -1.16%  "

$ws.Cells.Item(41, 2).Value = "MXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.416"
$ws.Cells.Item(41, 5).Value = "  -0.81%  "

$ws.Cells.Item(42, 2).Value = "PaxDollar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.003"
$ws.Cells.Item(42, 5).Value = "  +0.51%  "

$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.01503"
$ws.Cells.Item(43, 5).Value = "  +3.52%  "

$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.667"
$ws.Cells.Item(44, 5).Value = "  +11.68%  "

$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.3906"
$ws.Cells.Item(45, 5).Value = "  +4.75%  "

$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "6.941"
$ws.Cells.Item(46, 5).Value = "  +12.81%  "

$ws.Cells.Item(47, 2).Value = "Algorand"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1193"
$ws.Cells.Item(47, 5).Value = "  +6.78%  "

$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05330"
$ws.Cells.Item(48, 5).Value = "  +0.70%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.892"
$ws.Cells.Item(49, 5).Value = "  +3.43%  "

$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "30.79"
$ws.Cells.Item(50, 5).Value = "  -1.84%  "

$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.273"
$ws.Cells.Item(51, 5).Value = "  +5.10%  "
